$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update I-column (High_Prob_Choice_Corr reaction-time-like values) for rows 2-103 ---
$ws.Range("I2").Value = 4.4700000000000006
$ws.Range("I3").Value = 3.73
$ws.Range("I4").Value = 3.52
$ws.Range("I5").Value = 5.01
$ws.Range("I6").Value = 4.2200000000000006
$ws.Range("I7").Value = 3.57
$ws.Range("I8").Value = 3.97
$ws.Range("I9").Value = 4.1099999999999994
$ws.Range("I10").Value = 4.3499999999999996
$ws.Range("I11").Value = 3.7
$ws.Range("I12").Value = 3.8
$ws.Range("I13").Value = 5.01
$ws.Range("I14").Value = 4.25
$ws.Range("I15").Value = 4.1500000000000004
$ws.Range("I16").Value = 5.01
$ws.Range("I17").Value = 4.0299999999999994
$ws.Range("I18").Value = 3.95
$ws.Range("I19").Value = 4.57
$ws.Range("I20").Value = 4.4700000000000006
$ws.Range("I21").Value = 4.1099999999999994
$ws.Range("I22").Value = 5.01
$ws.Range("I23").Value = 4.7200000000000006
$ws.Range("I24").Value = 3.6
$ws.Range("I25").Value = 4.1400000000000006
$ws.Range("I26").Value = 5.01
$ws.Range("I27").Value = 3.61
$ws.Range("I28").Value = 5.01
$ws.Range("I29").Value = 4.43
$ws.Range("I30").Value = 4.1099999999999994
$ws.Range("I31").Value = 3.86
$ws.Range("I32").Value = 4.76
$ws.Range("I33").Value = 4.55
$ws.Range("I34").Value = 4.46
$ws.Range("I35").Value = 4.54
$ws.Range("I36").Value = 4.07
$ws.Range("I37").Value = 4.0199999999999996
$ws.Range("I38").Value = 3.92
$ws.Range("I39").Value = 4.7799999999999994
$ws.Range("I40").Value = 4.1400000000000006
$ws.Range("I41").Value = 3.67
$ws.Range("I42").Value = 4.42
$ws.Range("I43").Value = 5.01
$ws.Range("I44").Value = 3.89
$ws.Range("I45").Value = 4.08
$ws.Range("I46").Value = 4.04
$ws.Range("I47").Value = 4.13
$ws.Range("I48").Value = 4.29
$ws.Range("I49").Value = 4.59
$ws.Range("I50").Value = 4.3100000000000005
$ws.Range("I51").Value = 5.01
$ws.Range("I52").Value = 4.62
$ws.Range("I53").Value = 3.98
$ws.Range("I54").Value = 4.57
$ws.Range("I55").Value = 5.01
$ws.Range("I56").Value = 4.04
$ws.Range("I57").Value = 4.4800000000000004
$ws.Range("I58").Value = 3.74
$ws.Range("I59").Value = 4.42
$ws.Range("I60").Value = 4.09
$ws.Range("I61").Value = 5.01
$ws.Range("I62").Value = 5.01
$ws.Range("I63").Value = 4.29
$ws.Range("I64").Value = 4.7
$ws.Range("I65").Value = 3.57
$ws.Range("I66").Value = 3.86
$ws.Range("I67").Value = 5.01
$ws.Range("I68").Value = 3.66
$ws.Range("I69").Value = 4.21
$ws.Range("I70").Value = 5.01
$ws.Range("I71").Value = 4.38
$ws.Range("I72").Value = 4.38
$ws.Range("I73").Value = 3.59
$ws.Range("I74").Value = 3.95
$ws.Range("I75").Value = 3.9
$ws.Range("I76").Value = 4.32
$ws.Range("I77").Value = 4.38
$ws.Range("I78").Value = 4.42
$ws.Range("I79").Value = 3.64
$ws.Range("I80").Value = 4.1899999999999995
$ws.Range("I81").Value = 3.57
$ws.Range("I82").Value = 3.7
$ws.Range("I83").Value = 4.49
$ws.Range("I84").Value = 4
$ws.Range("I85").Value = 4.51
$ws.Range("I86").Value = 3.89
$ws.Range("I87").Value = 3.66
$ws.Range("I88").Value = 3.6
$ws.Range("I89").Value = 3.76
$ws.Range("I90").Value = 4.7799999999999994
$ws.Range("I91").Value = 4.76
$ws.Range("I92").Value = 4.37
$ws.Range("I93").Value = 3.6
$ws.Range("I94").Value = 4.1400000000000006
$ws.Range("I95").Value = 3.87
$ws.Range("I96").Value = 3.64
$ws.Range("I97").Value = 3.62
$ws.Range("I98").Value = 4.04
$ws.Range("I99").Value = 4.0299999999999994
$ws.Range("I100").Value = 4.5199999999999996
$ws.Range("I101").Value = 4.17
$ws.Range("I102").Value = 4.13
$ws.Range("I103").Value = 3.74

# --- Append 5 new trial rows (104-108) ---
$ws.Cells.Item(104,1).Value = 103
$ws.Cells.Item(104,2).Value = "Fractals/Version_2/ApoAV-250719-154.jpg"
$ws.Cells.Item(104,3).Value = "Reversal"
$ws.Cells.Item(104,4).Value = "75_Threat"
$ws.Cells.Item(104,5).Value = 0
$ws.Cells.Item(104,6).Value = 0
$ws.Cells.Item(104,7).Value = 0
$ws.Cells.Item(104,8).Value = "right"
$ws.Cells.Item(104,9).Value = 3.63
$ws.Cells.Item(104,10).Value = 0.99
$ws.Cells.Item(104,11).Value = 1

$ws.Cells.Item(105,1).Value = 104
$ws.Cells.Item(105,2).Value = "Fractals/Version_2/ApoAV-250719-72.jpg"
$ws.Cells.Item(105,3).Value = "Stable"
$ws.Cells.Item(105,4).Value = "75_Safe"
$ws.Cells.Item(105,5).Value = 0
$ws.Cells.Item(105,6).Value = 0
$ws.Cells.Item(105,7).Value = 0
$ws.Cells.Item(105,8).Value = "left"
$ws.Cells.Item(105,9).Value = 3.6
$ws.Cells.Item(105,10).Value = 1.43
$ws.Cells.Item(105,11).Value = 1

$ws.Cells.Item(106,1).Value = 105
$ws.Cells.Item(106,2).Value = "Fractals/Version_2/ApoAV-250719-154.jpg"
$ws.Cells.Item(106,3).Value = "Reversal"
$ws.Cells.Item(106,4).Value = "75_Threat"
$ws.Cells.Item(106,5).Value = 0
$ws.Cells.Item(106,6).Value = 0
$ws.Cells.Item(106,7).Value = 1
$ws.Cells.Item(106,8).Value = "right"
$ws.Cells.Item(106,9).Value = 4.17
$ws.Cells.Item(106,10).Value = 1.99
$ws.Cells.Item(106,11).Value = 1

$ws.Cells.Item(107,1).Value = 106
$ws.Cells.Item(107,2).Value = "Fractals/Version_2/ApoAV-250719-72.jpg"
$ws.Cells.Item(107,3).Value = "Stable"
$ws.Cells.Item(107,4).Value = "75_Safe"
$ws.Cells.Item(107,5).Value = 0
$ws.Cells.Item(107,6).Value = 0
$ws.Cells.Item(107,7).Value = 1
$ws.Cells.Item(107,8).Value = "left"
$ws.Cells.Item(107,9).Value = 4.73
$ws.Cells.Item(107,10).Value = 1.3699999999999999
$ws.Cells.Item(107,11).Value = 1

$ws.Cells.Item(108,1).Value = 107
$ws.Cells.Item(108,2).Value = "Fractals/Version_2/ApoAV-250719-154.jpg"
$ws.Cells.Item(108,3).Value = "Reversal"
$ws.Cells.Item(108,4).Value = "75_Threat"
$ws.Cells.Item(108,5).Value = 0
$ws.Cells.Item(108,6).Value = 0
$ws.Cells.Item(108,7).Value = 1
$ws.Cells.Item(108,8).Value = "right"
$ws.Cells.Item(108,9).Value = 3.64
$ws.Cells.Item(108,10).Value = 2.29
$ws.Cells.Item(108,11).Value = 1

# --- Update sheet view: scroll position + selection (matches saved view state) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 67
$win.ScrollColumn = 1
[void]$ws.Range("N104").Select()
